# Weekly update: insert two new price records (week of 2022-01-08, serial 44559)
# for "Betarraga" at "Femacal de La Calera", shifting all the existing rows
# below row 377 down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 378, pushing every row
# currently at 378 or below down by 2 (this is what creates the new
# dimension A1:R508 and shifts all subsequent data).
$ws.Rows("378:379").Insert()

# Fill in the first new row (378) - "Primera" quality entry for the new date.
$ws.Range("A378").Value = 3
$ws.Range("B378").Value = "Femacal de La Calera"
$ws.Range("C378").Value = "Coquimbo"
$ws.Range("D378").Value = 44559
$ws.Range("E378").Value = 5
$ws.Range("F378").Value = 100114014
$ws.Range("G378").Value = "Betarraga"
$ws.Range("H378").Value = "Sin especificar"
$ws.Range("I378").Value = "Primera"
$ws.Range("J378").Value = 3200
$ws.Range("K378").Value = 450
$ws.Range("L378").Value = 500
$ws.Range("M378").Value = 473
$ws.Range("N378").Value = "`$/paquete 4 unidades"
$ws.Range("O378").Value = "Provincia de Quillota"
$ws.Range("P378").Value = 118
$ws.Range("Q378").Value = 4
$ws.Range("R378").Value = "Hortaliza"

# Fill in the second new row (379) - "Segunda" quality entry for the new date.
$ws.Range("A379").Value = 3
$ws.Range("B379").Value = "Femacal de La Calera"
$ws.Range("C379").Value = "Coquimbo"
$ws.Range("D379").Value = 44559
$ws.Range("E379").Value = 5
$ws.Range("F379").Value = 100114014
$ws.Range("G379").Value = "Betarraga"
$ws.Range("H379").Value = "Sin especificar"
$ws.Range("I379").Value = "Segunda"
$ws.Range("J379").Value = 2800
$ws.Range("K379").Value = 300
$ws.Range("L379").Value = 350
$ws.Range("M379").Value = 332
$ws.Range("N379").Value = "`$/paquete 4 unidades"
$ws.Range("O379").Value = "Provincia de Quillota"
$ws.Range("P379").Value = 83
$ws.Range("Q379").Value = 4
$ws.Range("R379").Value = "Hortaliza"
